$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) for rows 2-8: 45185 -> 45204
$ws.Range("C2:C8").Value = 45204
